$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.270586252212524
$ws.Range("B1").Value = 2.560418844223022
$ws.Range("C1").Value = 2.200759649276733
$ws.Range("D1").Value = 2.334538459777832
$ws.Range("E1").Value = 2.784363508224487
